$p = $ppt.ActivePresentation

# --- Slide 10: "Ket qua" slide -----------------------------------------
# TextBox 14 holds "Do chinh xac : 58,42%" split across many runs; the
# authored edit retypes the "58,42%" portion so PowerPoint re-merges it
# into a single run (keeping the bold formatting of that portion).
$s10 = $p.Slides.Item(10)
$accShape = $s10.Shapes.Item(3)
$accRange = $accShape.TextFrame.TextRange
$accText = $accRange.Text
$accTarget = "58,42%"
$accStart = $accText.IndexOf($accTarget) + 1
$accRange.Characters($accStart, $accTarget.Length).Text = $accTarget

# --- Slide 2: "Bai toan" slide ------------------------------------------
# Content Placeholder 2 holds the long prompt sentence. The authored edit
# retypes two spans of that sentence (" (Gold ETF) 29 " and the
# "...gia vang" tail), which makes PowerPoint re-merge those runs while
# leaving the rest of the sentence's run boundaries untouched.
$s2 = $p.Slides.Item(2)
$bodyShape = $s2.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyText = $bodyRange.Text

$span1 = " (Gold ETF) 29 "
$start1 = $bodyText.IndexOf($span1) + 1
$bodyRange.Characters($start1, $span1.Length).Text = $span1

$span2 = "àm sao dự đoán được giá vàng"
$start2 = $bodyText.IndexOf($span2) + 1
$bodyRange.Characters($start2, $span2.Length).Text = $span2
